$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 22:16"

$data = @(
  ,("China",81008,41,71740,6013,1927,7,3255)
  ,("Italia",53578,6557,6072,42681,2857,793,4825)
  ,("España",25374,3803,2125,21871,1612,285,1378)
  ,("Estados Unidos",24142,4759,171,23683,64,32,288)
  ,("Alemania",22213,2365,209,21920,2,16,84)
  ,("Iran",20610,966,7635,11419,0,123,1556)
  ,("Francia",14459,1847,1587,12310,1525,112,562)
  ,("Corea del Sur",8799,147,2612,6085,59,8,102)
  ,("Suiza",6746,1131,15,6652,141,23,79)
  ,("Reino Unido",5018,1035,93,4692,20,56,233)
  ,("Paises Bajos",3631,637,2,3493,354,30,136)
  ,("Austria",2955,306,9,2938,15,2,8)
  ,("Belgica",2815,558,263,2485,288,30,67)
  ,("Noruega",2150,191,6,2137,28,0,7)
  ,("Suecia",1770,131,16,1734,71,4,20)
  ,("Dinamarca",1326,71,1,1312,42,4,13)
  ,("Portugal",1280,260,5,1263,26,6,12)
  ,("Canada",1280,193,14,1248,1,6,18)
  ,("Malasia",1183,153,114,1061,26,5,8)
  ,("Brasil",1128,158,2,1108,18,7,18)
  ,("Australia",1072,144,46,1019,2,0,7)
  ,("Japon",1046,39,215,795,55,1,36)
  ,("Chequia",995,162,6,989,7,0,0)
  ,("Turquia",947,277,0,926,0,12,21)
  ,("Israel",883,178,36,846,15,0,1)
  ,("Irlanda",785,102,5,777,13,0,3)
  ,("Pakistan",734,233,13,718,0,0,3)
  ,("Crucero",712,0,567,137,15,0,8)
  ,("Luxemburgo",670,186,6,656,3,3,8)
  ,("Chile",537,103,8,528,7,1,1)
  ,("Polonia",536,111,13,518,3,0,5)
  ,("Grecia",530,35,19,498,20,3,13)
  ,("Finlandia",523,73,10,512,2,1,1)
  ,("Ecuador",506,80,3,496,2,0,7)
  ,("Catar",481,11,27,454,6,0,0)
  ,("Islandia",473,64,5,467,1,1,1)
  ,("Indonesia",450,81,20,392,0,6,38)
  ,("Singapur",432,47,140,290,14,2,2)
  ,("Tailandia",411,89,44,366,7,0,1)
  ,("Arabia Saudita",392,48,16,376,0,0,0)
  ,("Eslovenia",383,42,0,382,12,0,1)
  ,("Rumania",367,59,52,315,14,0,0)
  ,("India",332,83,23,304,0,0,5)
  ,("Peru",318,55,1,312,5,1,5)
  ,("Barein",310,12,125,184,4,0,1)
  ,("Filipinas",307,77,13,275,1,1,19)
  ,("Estonia",306,23,2,304,0,0,0)
  ,("Rusia",306,53,16,289,0,0,1)
  ,("Egipto",294,9,42,242,0,2,10)
  ,("Hong Kong",273,17,98,171,4,0,4)
  ,("Sudafrica",240,38,2,238,0,0,0)
  ,("Libano",230,53,8,218,4,0,4)
  ,("Irak",214,6,51,146,0,0,17)
  ,("Croacia",206,76,5,200,0,0,1)
  ,("Mexico",203,39,4,197,1,1,2)
  ,("Panama",200,0,1,198,7,0,1)
  ,("Colombia",196,51,1,195,0,0,0)
  ,("Eslovaquia",178,41,7,171,2,0,0)
  ,("Kuwait",176,17,27,149,5,0,0)
  ,("Serbia",171,36,2,168,4,0,1)
  ,("Bulgaria",163,36,3,157,3,0,3)
  ,("Armenia",160,24,1,159,2,0,0)
  ,("San Marino",160,9,4,136,12,6,20)
  ,("Argentina",158,0,3,151,0,1,4)
  ,("Taiwan",153,18,28,123,0,0,2)
  ,("Emiratos Arabes Unidos",153,13,38,113,2,0,2)
  ,("Argelia",139,45,65,59,0,4,15)
  ,("Letonia",124,13,1,123,0,0,0)
  ,("Costa Rica",117,4,2,113,2,0,2)
  ,("Republica Dominicana",112,40,0,109,0,1,3)
  ,("Uruguay",110,0,0,110,0,0,0)
  ,("Hungria",103,18,7,92,6,0,4)
  ,("Jordania",99,15,1,98,0,0,0)
  ,("Lituania",99,30,1,97,1,0,1)
  ,("Marruecos",96,10,3,90,1,0,3)
  ,("Vietnam",94,3,17,77,2,0,0)
  ,("Bosnia y Herzegovina",93,4,2,90,1,1,1)
  ,("Islas Feroe",92,12,3,89,0,0,0)
  ,("Principado de Andorra",88,13,1,87,2,0,0)
  ,("Republica de Macedonia",85,9,1,84,1,0,0)
  ,("Republica de Chipre",84,9,3,81,3,0,0)
  ,("Brunei",83,5,1,82,2,0,0)
  ,("Moldavia",80,14,1,78,3,0,1)
  ,("Sri Lanka",77,4,3,74,2,0,0)
  ,("Albania",76,6,2,72,2,0,2)
  ,("Bielorrusia",76,7,15,61,0,0,0)
  ,("Malta",73,9,2,71,1,0,0)
  ,("Venezuela",70,5,15,55,2,0,0)
  ,("Burkina Faso",64,24,5,56,0,2,3)
  ,("Tunez",60,6,1,58,7,0,1)
  ,("Senegal",56,9,5,51,0,0,0)
  ,("Kazajistan",54,2,0,54,0,0,0)
  ,("Camboya",53,2,2,51,0,0,0)
  ,("Azerbaiyan",53,9,11,41,0,0,1)
  ,("Estado de Palestina",53,5,17,36,0,0,0)
  ,("Nueva Zelanda",52,0,0,52,0,0,0)
  ,("Oman",52,4,13,39,0,0,0)
  ,("Guadalupe",51,0,0,50,4,0,1)
  ,("Trinidad yTobago",49,40,0,49,0,0,0)
  ,("Georgia",49,5,1,48,1,0,0)
  ,("Ucrania",47,6,1,43,0,0,3)
  ,("Reunion",45,7,0,45,0,0,0)
  ,("Uzbekistan",41,8,0,41,0,0,0)
  ,("Liechtenstein",37,9,0,37,0,0,0)
  ,("Martinica",32,0,0,31,7,0,1)
  ,("Camerun",27,0,2,25,0,0,0)
  ,("Honduras",24,0,0,24,0,0,0)
  ,("Afganistan",24,0,1,23,0,0,0)
  ,("Banglades",24,4,3,19,0,1,2)
  ,("Consejo Danes para los Refugiados",23,5,0,22,0,1,1)
  ,("Nigeria",22,10,1,21,0,0,0)
  ,("Puerto Rico",21,7,0,21,0,0,0)
  ,("Cuba",21,0,0,20,0,0,1)
  ,("Bolivia",19,3,0,19,0,0,0)
  ,("Ghana",19,3,0,18,0,1,1)
  ,("Jamaica",19,0,2,16,0,0,1)
  ,("Guayana Francesa",18,3,0,18,0,0,0)
  ,("Monaco",18,7,1,17,0,0,0)
  ,("Paraguay",18,0,0,17,1,1,1)
  ,("Macao",18,1,10,8,0,0,0)
  ,("Ruanda",17,0,0,17,0,0,0)
  ,("Guatemala",17,5,0,16,0,0,1)
  ,("Togo",16,7,0,16,0,0,0)
  ,("Montenegro",16,2,0,16,0,0,0)
  ,("Guam",15,1,0,15,0,0,0)
  ,("Polinesia Francesa",15,4,0,15,0,0,0)
  ,("Kirguistan",14,8,0,14,0,0,0)
  ,("Mauricio",14,2,0,13,0,1,1)
  ,("Costa de Marfil",14,5,1,13,0,0,0)
  ,("Maldivas",13,0,3,10,0,0,0)
  ,("Mongolia",10,4,0,10,0,0,0)
  ,("Gibraltar",10,0,2,8,0,0,0)
  ,("Etiopia",9,0,0,9,0,0,0)
  ,("Seychelles",7,0,0,7,0,0,0)
  ,("Mayotte",7,0,0,7,0,0,0)
  ,("Kenia",7,0,0,7,0,0,0)
  ,("Barbados",6,0,0,6,0,0,0)
  ,("Guinea Ecuatorial",6,0,0,6,0,0,0)
  ,("Islas Virgenes de los Estados Unidos",6,3,0,6,0,0,0)
  ,("Tanzania",6,0,0,6,0,0,0)
  ,("Surinam",5,1,0,5,0,0,0)
  ,("Gabon",5,1,0,4,0,0,1)
  ,("Aruba",5,0,1,4,0,0,0)
  ,("Guyana",5,0,0,4,0,0,1)
  ,("Bahamas",4,0,0,4,0,0,0)
  ,("Nueva Caledonia",4,2,0,4,0,0,0)
  ,("San Martin (Parte Francesa)",4,0,0,4,0,0,0)
  ,("Cabo Verde",3,2,0,3,0,0,0)
  ,("Zimbabue",3,2,0,3,0,0,0)
  ,("Namibia",3,0,0,3,0,0,0)
  ,("Madagascar",3,0,0,3,0,0,0)
  ,("Congo",3,0,0,3,0,0,0)
  ,("Liberia",3,1,0,3,0,0,0)
  ,("El Salvador",3,2,0,3,0,0,0)
  ,("Republica de Africa Central",3,0,0,3,0,0,0)
  ,("San Bartolome",3,0,0,3,0,0,0)
  ,("Curazao",3,0,0,2,0,0,1)
  ,("Islas Caimanes",3,0,0,2,0,0,1)
  ,("Haiti",2,0,0,2,0,0,0)
  ,("Groenlandia",2,0,0,2,0,0,0)
  ,("Mauritania",2,0,0,2,0,0,0)
  ,("Santa Lucia",2,0,0,2,0,0,0)
  ,("Benin",2,0,0,2,0,0,0)
  ,("Isla de Man",2,0,0,2,0,0,0)
  ,("Butan",2,0,0,2,0,0,0)
  ,("Guinea",2,0,0,2,0,0,0)
  ,("Bermudas",2,0,0,2,0,0,0)
  ,("Nicaragua",2,0,0,2,0,0,0)
  ,("Zambia",2,0,0,2,0,0,0)
  ,("Fiyi",2,1,0,2,0,0,0)
  ,("Angola",2,1,0,2,0,0,0)
  ,("Sudan",2,0,0,1,0,0,1)
  ,("San Vicente y las Granadinas",1,0,0,1,0,0,0)
  ,("Santa Sede",1,0,0,1,0,0,0)
  ,("Antigua y Barbuda",1,0,0,1,0,0,0)
  ,("Gambia",1,0,0,1,0,0,0)
  ,("Montserrat",1,0,0,1,0,0,0)
  ,("Suazilandia",1,0,0,1,0,0,0)
  ,("San Martin (Parte Holandesa)",1,0,0,1,0,0,0)
  ,("Uganda",1,1,0,1,0,0,0)
  ,("Republica del Chad",1,0,0,1,0,0,0)
  ,("Niger",1,0,0,1,0,0,0)
  ,("Eritrea",1,1,0,1,0,0,0)
  ,("Republica de Yibuti",1,0,0,1,0,0,0)
  ,("Papua Nueva Guinea",1,0,0,1,0,0,0)
  ,("Somalia",1,0,0,1,0,0,0)
  ,("Timor Oriental",1,1,0,1,0,0,0)
  ,("Nepal",1,0,1,0,0,0,0)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
    $ws.Cells.Item($row, 8).Value = $item[7]
    $row++
}
